$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (matches the existing D1:G1 style - numFmt 164, default font)
$ws.Range("H1").Value = "context"
$ws.Range("I1").Value = "type"
$ws.Range("H1:I1").NumberFormat = "0.0000"

# Hill II / mudbrick rows (AA-1, AA-2, AA-3)
$ws.Range("H2").Value = "Hill II"
$ws.Range("I2").Value = "mudbrick"
$ws.Range("H3").Value = "Hill II"
$ws.Range("I3").Value = "mudbrick"
$ws.Range("H4").Value = "Hill II"
$ws.Range("I4").Value = "mudbrick"

# Urartian silo / mudbrick rows (AA-4, AA-4 (2))
$ws.Range("H6").Value = "Urartian silo"
$ws.Range("I6").Value = "mudbrick"
$ws.Range("H7").Value = "Urartian silo"
$ws.Range("I7").Value = "mudbrick"

# Phase I / mudbrick rows (AA-5, AA-6)
$ws.Range("H8").Value = "Phase I"
$ws.Range("I8").Value = "mudbrick"
$ws.Range("H9").Value = "Phase I"
$ws.Range("I9").Value = "mudbrick"

# Phase II or III / mudbrick rows (AA-7, AA-8)
$ws.Range("H10").Value = "Phase II or III"
$ws.Range("I10").Value = "mudbrick"
$ws.Range("H11").Value = "Phase II or III"
$ws.Range("I11").Value = "mudbrick"

# Phase II / mudbrick rows (AA-9, AA-10)
$ws.Range("H12").Value = "Phase II"
$ws.Range("I12").Value = "mudbrick"
$ws.Range("H13").Value = "Phase II"
$ws.Range("I13").Value = "mudbrick"

# kekkila control row gets NA/NA last (so "NA" is appended last to sharedStrings)
$ws.Range("H5").Value = "NA"
$ws.Range("I5").Value = "NA"

# Re-apply font on the new cells (mirrors the author re-typing the Normal font,
# which Excel records as a fresh font entry + new cellXfs records)
$ws.Range("H2:I4").Font.Name = "Calibri"
$ws.Range("H6:I13").Font.Name = "Calibri"
$ws.Range("H5:I5").NumberFormat = "0.0000"
$ws.Range("H5:I5").Font.Name = "Calibri"

# Fit the new context column to its content
$ws.Columns("H:H").AutoFit() | Out-Null

# Match the author's final selection
$ws.Range("I6").Select() | Out-Null
